# feat: add 2022-Q4 data
#
# Before: 总计, 2022-Q3, 2022-Q2
# After:  总计, 2022-Q4, 2022-Q3, 2022-Q2
#
# The existing "2022-Q3" fund-holdings sheet is duplicated (to inherit its
# header/formatting) into the new slot right before it, renamed to
# "2022-Q4", and its values overwritten with the new quarter's numbers.
# The original "2022-Q3" sheet is left untouched (just shifted one tab to
# the right, same as "2022-Q2"). The "总计" roll-up sheet gets a new row
# for the Q3 total (which used to live in row 2) and row 2 itself becomes
# the Q4 total; the Q2 row shifts down to row 4.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q4" sheet by copying "2022-Q3" so it keeps the
#    identical header/borders/font styling, placed immediately before the
#    existing "2022-Q3" tab.
# ---------------------------------------------------------------------
$wsQ3 = $wb.Worksheets.Item("2022-Q3")
$wsQ3.Copy($wsQ3, $null)
$wsQ4 = $wb.Worksheets.Item(2)
$wsQ4.Name = "2022-Q4"

# Fund-code / name / ratio columns are stored as text in this workbook
# (so things like leading zeros in fund codes and fixed-decimal ratios
# survive round-tripping) - force text formatting before writing so Excel
# doesn't silently re-interpret them as numbers.
$wsQ4.Range("B2:G6").NumberFormat = "@"

# Row 2: 009010 华夏兴阳一年持有期混合
$wsQ4.Range("B2").Value = "009010"
$wsQ4.Range("C2").Value = "华夏兴阳一年持有期混合"
$wsQ4.Range("D2").Value = "27.18"
$wsQ4.Range("E2").Value = "90.65"
$wsQ4.Range("F2").Value = "2.23"
$wsQ4.Range("G2").Value = "0.6061"
$wsQ4.Range("H2").Value = 10

# Row 3: 160322 华夏港股通精选股票（LOF）A
$wsQ4.Range("B3").Value = "160322"
$wsQ4.Range("C3").Value = "华夏港股通精选股票（LOF）A"
$wsQ4.Range("D3").Value = "13.82"
$wsQ4.Range("E3").Value = "92.59"
$wsQ4.Range("F3").Value = "2.68"
$wsQ4.Range("G3").Value = "0.3704"
$wsQ4.Range("H3").Value = 6

# Row 4: 012884 华夏港股通精选股票（LOF）C
$wsQ4.Range("B4").Value = "012884"
$wsQ4.Range("C4").Value = "华夏港股通精选股票（LOF）C"
$wsQ4.Range("D4").Value = "0.69"
$wsQ4.Range("E4").Value = "92.59"
$wsQ4.Range("F4").Value = "2.68"
$wsQ4.Range("G4").Value = "0.0185"
$wsQ4.Range("H4").Value = 6

# Row 5: 005255 浦银安盛港股通量化混合A
$wsQ4.Range("B5").Value = "005255"
$wsQ4.Range("C5").Value = "浦银安盛港股通量化混合A"
$wsQ4.Range("D5").Value = "0.34"
$wsQ4.Range("E5").Value = "59.70"
$wsQ4.Range("F5").Value = "2.71"
$wsQ4.Range("G5").Value = "0.0092"
$wsQ4.Range("H5").Value = 7

# Row 6: 013224 浦银安盛港股通量化混合C
$wsQ4.Range("B6").Value = "013224"
$wsQ4.Range("C6").Value = "浦银安盛港股通量化混合C"
$wsQ4.Range("D6").Value = "0.17"
$wsQ4.Range("E6").Value = "59.70"
$wsQ4.Range("F6").Value = "2.71"
$wsQ4.Range("G6").Value = "0.0046"
$wsQ4.Range("H6").Value = 7

# ---------------------------------------------------------------------
# 2. Update the "总计" roll-up sheet: keep row 1 headers, turn row 2 into
#    the Q4 totals, insert a Q3 row (the old row-2 values) and push the
#    Q2 row down to row 4.
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

# Copy column-A's number style down onto the two rows we are about to
# populate so they keep the same (bordered/centered) look as row 2.
$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3:A4").PasteSpecial(-4122)

# Row 4 becomes the old Q2 row (shifted down from row 3).
$wsTotal.Range("A4").Value = 2
$wsTotal.Range("B4").Value = "2022-Q2"
$wsTotal.Range("C4").Value = 2
$wsTotal.Range("D4").Value = 0.03

# Row 3 becomes the old Q3 row (shifted down from row 2).
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q3"
$wsTotal.Range("C3").Value = 5
$wsTotal.Range("D3").Value = 0.85

# Row 2 becomes the new Q4 row.
$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 5
$wsTotal.Range("D2").Value = 1.01
